# Apply corrected ordering of the "Recorded By" (column G) values.
# This reorders the comma-separated list of recorders for the affected rows
# (the underlying "System" entry moves to the front of the list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = "System, backup@backdoor.com, system"
$ws.Range("G3").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G4").Value2 = "System, backup@backdoor.com"
$ws.Range("G6").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G10").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G11").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G12").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G13").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G14").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G15").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G17").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G18").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G19").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G20").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G21").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G22").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G24").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G29").Value2 = "System, backup@backdoor.com, system"
$ws.Range("G30").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G31").Value2 = "System, backup@backdoor.com"
$ws.Range("G33").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G37").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G38").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G39").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G40").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G41").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G42").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G44").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G45").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G46").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G47").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G48").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G49").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G51").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G56").Value2 = "System, backup@backdoor.com, system"
$ws.Range("G57").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G58").Value2 = "System, backup@backdoor.com"
$ws.Range("G60").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G64").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G65").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G66").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G67").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G68").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G69").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G71").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G72").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G73").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G74").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G75").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G76").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G78").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G86").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G87").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G88").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G89").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G93").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G95").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G96").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G97").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G99").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G102").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G112").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G113").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G114").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G115").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G119").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G121").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G122").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G123").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G125").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G128").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G138").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G139").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G140").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G141").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G145").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G147").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G148").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G149").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G151").Value2 = "System, dnasr281@gmail.com"
$ws.Range("G154").Value2 = "System, dnasr281@gmail.com"
